$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.838.45'
$ws.Range('E2').Value = '  +8.35%  '

$ws.Range('D3').Value = '3.222.58'
$ws.Range('E3').Value = '  +3.88%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '''400.08'
$ws.Range('E5').Value = '  +3.56%  '

$ws.Range('D6').Value = '''110.42'
$ws.Range('E6').Value = '  +6.62%  '

$ws.Range('E7').Value = '  +2.61%  '

$ws.Range('E8').Value = '  -0.10%  '

$ws.Range('D9').Value = '''0.624'
$ws.Range('E9').Value = '  +6.50%  '

$ws.Range('D10').Value = '''39.36'
$ws.Range('E10').Value = '  +6.21%  '

$ws.Range('E11').Value = '  +4.84%  '

$ws.Range('E12').Value = '  +1.87%  '

$ws.Range('D13').Value = '3.734.12'
$ws.Range('E13').Value = '  +3.75%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '''19.11'
$ws.Range('E14').Value = '  +2.62%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '''8.08'
$ws.Range('E15').Value = '  +2.91%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''1.06'
$ws.Range('E16').Value = '  +6.48%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.219.44'
$ws.Range('E17').Value = '  +3.90%  '

$ws.Range('D18').Value = '''10.69'
$ws.Range('E18').Value = '  -3.15%  '

$ws.Range('D19').Value = '55.715.82'
$ws.Range('E19').Value = '  +7.99%  '

$ws.Range('D20').Value = '''3.40'
$ws.Range('E20').Value = '  +2.80%  '

$ws.Range('D21').Value = '''0.0000103'
$ws.Range('E21').Value = '  +6.31%  '

$ws.Range('D22').Value = '''13.05'
$ws.Range('E22').Value = '  +5.16%  '

$ws.Range('D23').Value = '''300.73'
$ws.Range('E23').Value = '  +12.92%  '

$ws.Range('D24').Value = '''74.71'
$ws.Range('E24').Value = '  +6.67%  '

$ws.Range('E25').Value = '  +2.38%  '

$ws.Range('D26').Value = '''8.14'
$ws.Range('E26').Value = '  +0.22%  '

$ws.Range('D27').Value = '''28.41'
$ws.Range('E27').Value = '  +4.77%  '

$ws.Range('D28').Value = '''7.44'
$ws.Range('E28').Value = '  +2.25%  '

$ws.Range('D29').Value = '''0.171'
$ws.Range('E29').Value = '  +1.78%  '

$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '''11.44'
$ws.Range('E31').Value = '  +10.48%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.111'
$ws.Range('E32').Value = '  +4.01%  '

$ws.Range('D33').Value = '''0.0493'
$ws.Range('E33').Value = '  +5.58%  '

$ws.Range('D34').Value = '''36.51'
$ws.Range('E34').Value = '  +2.99%  '

$ws.Range('E35').Value = '  +2.19%  '

$ws.Range('D36').Value = '''51.31'
$ws.Range('E36').Value = '  +2.23%  '

$ws.Range('E37').Value = '  +5.49%  '

$ws.Range('E38').Value = '  +0.05%  '

$ws.Range('E39').Value = '  +22.39%  '

$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = '''4.08'
$ws.Range('E40').Value = '  +10.94%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '''133.25'
$ws.Range('E41').Value = '  +3.23%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''1.92'
$ws.Range('E42').Value = '  +1.61%  '

$ws.Range('D43').Value = '''17.09'
$ws.Range('E43').Value = '  +3.19%  '

$ws.Range('E44').Value = '  +2.88%  '

$ws.Range('D45').Value = '''0.282'
$ws.Range('E45').Value = '  -4.70%  '

$ws.Range('D46').Value = '''22.38'
$ws.Range('E46').Value = '  -0.37%  '

$ws.Range('D47').Value = '''2.13'
$ws.Range('E47').Value = '  +43.17%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '''2.12'
$ws.Range('E48').Value = '  +2.26%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.157.43'
$ws.Range('E49').Value = '  +4.47%  '

$ws.Range('E50').Value = '  -0.17%  '

$ws.Range('E51').Value = '  +8.63%  '
